# Apply "correct weight optimization" edit to the Pull Sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for rows 2-40 (columns A-F)
# A = Pull #, B = Cable Size, C = Local / Express, D = From, E = To, F = Distance
$data = @(
    @(1,  "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(2,  "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(3,  "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(4,  "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(5,  "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(6,  "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(7,  "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(8,  "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(9,  "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(10, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(11, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(12, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(13, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(14, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(15, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(16, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(17, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(18, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(19, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(20, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(22, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(23, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(24, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(25, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(26, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(27, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(28, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(29, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(30, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(31, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(32, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(33, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(34, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(35, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(36, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(37, "STAR QUAD", "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(38, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(39, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300),
    @(40, "7C#14",     "EXPRESS", "SWITCH-HTR-A", "RELAY-RM-1", 300)
)

$rowNum = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowNum, 1).Value = $row[0]
    $ws.Cells.Item($rowNum, 2).Value = $row[1]
    $ws.Cells.Item($rowNum, 3).Value = $row[2]
    $ws.Cells.Item($rowNum, 4).Value = $row[3]
    $ws.Cells.Item($rowNum, 5).Value = $row[4]
    $ws.Cells.Item($rowNum, 6).Value = $row[5]
    $rowNum++
}

# Update the view: scroll position and selection, matching the author's saved state
$ws.Range("A14").Select()
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("C27").Select()
